$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 60/61 (pushes the existing rows 60-88 down to 62-90,
# matching the weekly-refresh pattern: a new pair of "Primera"/"Segunda"
# records for a newer date is prepended and the oldest rows survive,
# shifted, at the bottom).
$ws.Rows.Item(60).Resize(2).Insert()

# New row 60: Primera, fecha 2022-03-07 (serial 44627)
$ws.Cells.Item(60, 1).Value2 = 9
$ws.Cells.Item(60, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(60, 3).Value2 = "Metropolitana"
$ws.Cells.Item(60, 4).Value2 = 44627
$ws.Cells.Item(60, 5).Value2 = 13
$ws.Cells.Item(60, 6).Value2 = 100114007
$ws.Cells.Item(60, 7).Value2 = "Jengibre"
$ws.Cells.Item(60, 8).Value2 = "Sin especificar"
$ws.Cells.Item(60, 9).Value2 = "Primera"
$ws.Cells.Item(60, 10).Value2 = 790
$ws.Cells.Item(60, 11).Value2 = 14000
$ws.Cells.Item(60, 12).Value2 = 15000
$ws.Cells.Item(60, 13).Value2 = 14494
$ws.Cells.Item(60, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(60, 15).Value2 = "Perú"
$ws.Cells.Item(60, 16).Value2 = 1115
$ws.Cells.Item(60, 17).Value2 = 13
$ws.Cells.Item(60, 18).Value2 = "Hortaliza"

# New row 61: Segunda, fecha 2022-03-07 (serial 44627)
$ws.Cells.Item(61, 1).Value2 = 9
$ws.Cells.Item(61, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(61, 3).Value2 = "Metropolitana"
$ws.Cells.Item(61, 4).Value2 = 44627
$ws.Cells.Item(61, 5).Value2 = 13
$ws.Cells.Item(61, 6).Value2 = 100114007
$ws.Cells.Item(61, 7).Value2 = "Jengibre"
$ws.Cells.Item(61, 8).Value2 = "Sin especificar"
$ws.Cells.Item(61, 9).Value2 = "Segunda"
$ws.Cells.Item(61, 10).Value2 = 340
$ws.Cells.Item(61, 11).Value2 = 13000
$ws.Cells.Item(61, 12).Value2 = 13000
$ws.Cells.Item(61, 13).Value2 = 13000
$ws.Cells.Item(61, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(61, 15).Value2 = "Perú"
$ws.Cells.Item(61, 16).Value2 = 1000
$ws.Cells.Item(61, 17).Value2 = 13
$ws.Cells.Item(61, 18).Value2 = "Hortaliza"

Write-Output "done"
